$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update existing sound entry (heart beat.mp3) ---------------
$ws.Range("A2").Value = "heart beat.mp3"
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = "low-frequency beat heart heartbeat"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "This was produced by tapping on a stethoscope which had an earbud pressed against a Shure SM57 mic. Low pass filter applied, as well as compression and a gate. Chorus added. Used a recording cassette deck as a preamp going into an M-Audio Audiophile USB soundcard.`nNote: If you're having problems listening to this clip, the cutoff frequency of your speaker set may be too high(solution: new speakers). The signal strength exists almost entirely in the very low frequencies, so you may need a sub-woofer to hear it. Otherwise, try turning your speaker volume all the way up. Doing so may saturate the signal and at least allow you to hear the harmonics of the signal caused by the distortion. I don't recommend it, but you'll at least maybe be able to hear something."
$ws.Range("F2").Value = "Creative Commons 0"
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = 0

# --- Row 3: new sound entry (rbh thunder storm.wav) ---------------------
$ws.Range("A3").Value = "rbh thunder storm.wav"
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = "siren xy stereo heavy-rain thunder-storm storm nature field-recording rain thunder purist weather lightning"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "rain and several thunder claps. Stereo. Recorded with a minidisc recorder and an AT822 mic. Edited in Protools and saved as a 48/16 wav file.`nI have more sounds available here:`nhttp://sfx.TakomaMedia.com"
$ws.Range("F3").Value = "Attribution"
$ws.Range("G3").Value = "Storm"
$ws.Range("H3").Value = 0

# --- Row 4: new sound entry (STEREO_41.wav) ------------------------------
$ws.Range("A4").Value = "STEREO_41.wav"
$ws.Range("B4").Value = "Close seashore - Saintes-Maries-de-la-Mer"
$ws.Range("C4").Value = "shure-mv88 beach close waves sea seashore shore water field-recording"
$ws.Range("D4").Value = "43.44902, 4.40802, 16"
$ws.Range("E4").Value = "Close recording of the seashore in Saintes-Maries-de-la-Mer, Provence-Alpes-Côte d'Azur, France. Recorded with my mobile phone with a Shure mv88 on August 2015."
$ws.Range("F4").Value = "Attribution"
$ws.Range("G4").Value = "Provence-Alpes-Côte d'Azur"
$ws.Range("H4").Value = 0

# --- Row heights so the wrapped, multi-line text is fully visible -------
$ws.Rows.Item(2).RowHeight = 192
$ws.Rows.Item(3).RowHeight = 128
$ws.Rows.Item(4).RowHeight = 96

# --- Column E is now wide and wraps its text -----------------------------
$ws.Columns.Item(5).ColumnWidth = 69.1640625
$ws.Range("A1:H4").WrapText = $true

# --- Hyperlinks on the sound file name / name cells ----------------------
# Adding a hyperlink applies Excel's built-in "Hyperlink" style, so restore
# the original text + the sheet's normal (wrapped) formatting afterwards.
$ws.Hyperlinks.Add($ws.Range("B4"), "https://freesound.org/people/frederic.font/sounds/322271/", "", "", "https://freesound.org/people/frederic.font/sounds/322271/")
$ws.Range("B4").Value = "Close seashore - Saintes-Maries-de-la-Mer"
$ws.Range("F4").Copy()
$ws.Range("B4").PasteSpecial(-4122)

$ws.Hyperlinks.Add($ws.Range("A2"), "https://freesound.org/people/greyseraphim/sounds/21409/", "", "", "https://freesound.org/people/greyseraphim/sounds/21409/")
$ws.Range("A2").Value = "heart beat.mp3"
$ws.Range("F2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$excel.CutCopyMode = $false

$ws.Range("J2").Select()
